$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4042.261178533166
$ws.Range("C2").Value = 99.63138986942523
$ws.Range("D2").Value = 2872.758718594442
$ws.Range("B3").Value = 4843.104700189625
$ws.Range("C3").Value = 133.684749483326
$ws.Range("D3").Value = 3335.020037985748
$ws.Range("B4").Value = 5694.049009783843
$ws.Range("C4").Value = 132.0249906400697
$ws.Range("D4").Value = 3798.162259809709
$ws.Range("B5").Value = 6438.611309263215
$ws.Range("C5").Value = 146.4752784414153
$ws.Range("D5").Value = 4245.31533836317
$ws.Range("B6").Value = 7059.882565624657
$ws.Range("C6").Value = 163.6983367069226
$ws.Range("D6").Value = 4600.02708193398
$ws.Range("B7").Value = 7603.291433288357
$ws.Range("C7").Value = 173.9658407765299
$ws.Range("D7").Value = 4833.91336467765
$ws.Range("B8").Value = 8032.262111407125
$ws.Range("C8").Value = 190.6088938661355
$ws.Range("D8").Value = 5058.054074345984
$ws.Range("B9").Value = 8313.879475221485
$ws.Range("C9").Value = 199.6609295146157
$ws.Range("D9").Value = 5296.037405313211
$ws.Range("B10").Value = 8729.470243753516
$ws.Range("C10").Value = 212.6384564337307
$ws.Range("D10").Value = 5451.621370691195
$ws.Range("B11").Value = 9112.563602116979
$ws.Range("C11").Value = 238.6504178826896
$ws.Range("D11").Value = 5629.086028057161
$ws.Range("B12").Value = 9453.505841529721
$ws.Range("C12").Value = 243.6519269973406
$ws.Range("D12").Value = 5764.157092739942
$ws.Range("B13").Value = 9774.924001650776
$ws.Range("C13").Value = 258.6806115824122
$ws.Range("D13").Value = 5898.040155849277
$ws.Range("B14").Value = 10177.41970499595
$ws.Range("C14").Value = 277.3238092996165
$ws.Range("D14").Value = 6027.284754104194
$ws.Range("B15").Value = 10457.61164861421
$ws.Range("C15").Value = 285.959734488659
$ws.Range("D15").Value = 6129.933177734317
$ws.Range("B16").Value = 10739.56104438033
$ws.Range("C16").Value = 331.8596349857691
$ws.Range("D16").Value = 6242.907871355815
$ws.Range("B17").Value = 10991.37258145039
$ws.Range("C17").Value = 341.3306269331836
$ws.Range("D17").Value = 6313.858401919963
$ws.Range("B18").Value = 11177.68563349658
$ws.Range("C18").Value = 347.8315910842927
$ws.Range("D18").Value = 6338.408447373296
$ws.Range("B19").Value = 11417.61049189586
$ws.Range("C19").Value = 367.4942650344049
$ws.Range("D19").Value = 6366.887305438978
$ws.Range("B20").Value = 11595.06932656615
$ws.Range("C20").Value = 381.6807619050949
$ws.Range("D20").Value = 6446.558841296247
$ws.Range("B21").Value = 11777.512002425
$ws.Range("C21").Value = 384.4399606484644
$ws.Range("D21").Value = 6463.165291660873
$ws.Range("B22").Value = 11997.33539222781
$ws.Range("C22").Value = 391.6039521416105
$ws.Range("D22").Value = 6489.928969743357
$ws.Range("B23").Value = 12128.77514305258
$ws.Range("C23").Value = 398.3599774891508
$ws.Range("D23").Value = 6449.104745751848
$ws.Range("B24").Value = 12392.78709634662
$ws.Range("C24").Value = 414.1935364630534
$ws.Range("D24").Value = 6474.778370142047
$ws.Range("B25").Value = 12635.33807925911
$ws.Range("C25").Value = 427.5927052356462
$ws.Range("D25").Value = 6482.587705592893
$ws.Range("B26").Value = 12782.60586740629
$ws.Range("C26").Value = 435.5577939095213
$ws.Range("D26").Value = 6451.901730266419
$ws.Range("B27").Value = 12962.66049171196
$ws.Range("C27").Value = 447.4189084122374
$ws.Range("D27").Value = 6407.085948404506
$ws.Range("B28").Value = 13083.54188627047
$ws.Range("C28").Value = 461.2612545024355
$ws.Range("D28").Value = 6355.752448061746
$ws.Range("B29").Value = 13266.07278680368
$ws.Range("C29").Value = 468.3337443280442
$ws.Range("D29").Value = 6312.354616102235
$ws.Range("B30").Value = 13394.54002793477
$ws.Range("C30").Value = 478.5028027028271
$ws.Range("D30").Value = 6257.337492328858
$ws.Range("B31").Value = 13518.25868908472
$ws.Range("C31").Value = 485.0739095759201
$ws.Range("D31").Value = 6147.866861531261
$ws.Range("B32").Value = 13803.21419372251
$ws.Range("C32").Value = 459.04321535677
$ws.Range("D32").Value = 6122.895502285376
$ws.Range("B33").Value = 13952.14365383609
$ws.Range("C33").Value = 467.1473963087757
$ws.Range("D33").Value = 6049.27606144307
$ws.Range("B34").Value = 14030.55613453701
$ws.Range("C34").Value = 485.350364637565
$ws.Range("D34").Value = 5966.33096611939
$ws.Range("B35").Value = 14179.82315872458
$ws.Range("C35").Value = 505.1371702543942
$ws.Range("D35").Value = 5906.037691588646
$ws.Range("B36").Value = 14394.52610147837
$ws.Range("C36").Value = 506.8776412173643
$ws.Range("D36").Value = 5810.189512517667
$ws.Range("B37").Value = 14503.81357347808
$ws.Range("C37").Value = 520.8293881374873
$ws.Range("D37").Value = 5699.648857028375
$ws.Range("B38").Value = 14649.18554172886
$ws.Range("C38").Value = 528.5191406727265
$ws.Range("D38").Value = 5601.15213057723
$ws.Range("B39").Value = 14776.55981597375
$ws.Range("C39").Value = 540.6026932258874
$ws.Range("D39").Value = 5500.676729919274
$ws.Range("B40").Value = 14946.05459987146
$ws.Range("C40").Value = 554.2200523336803
$ws.Range("D40").Value = 5388.849305042385
$ws.Range("B41").Value = 15022.46331036672
$ws.Range("C41").Value = 589.6639227267999
$ws.Range("D41").Value = 5231.684012435935
$ws.Range("B42").Value = 15150.20247258959
$ws.Range("C42").Value = 599.6217302859592
$ws.Range("D42").Value = 5105.566710040019
$ws.Range("B43").Value = 15251.533915885
$ws.Range("C43").Value = 607.6474737616869
$ws.Range("D43").Value = 4967.486267451978
$ws.Range("B44").Value = 15353.17079266563
$ws.Range("C44").Value = 606.7788020492254
$ws.Range("D44").Value = 4824.690910951283
$ws.Range("B45").Value = 15495.65087828405
$ws.Range("C45").Value = 617.3690200466482
$ws.Range("D45").Value = 4674.281689593379
$ws.Range("B46").Value = 15547.2376000273
$ws.Range("C46").Value = 626.1883589175171
$ws.Range("D46").Value = 4712.605034101749
$ws.Range("B47").Value = 15724.50491018924
$ws.Range("C47").Value = 639.3888747846812
$ws.Range("D47").Value = 4564.993103103866
$ws.Range("B48").Value = 15825.63908986542
$ws.Range("C48").Value = 653.7152853326515
$ws.Range("D48").Value = 4359.242939525529
$ws.Range("B49").Value = 15816.64570403627
$ws.Range("C49").Value = 658.9732495405109
$ws.Range("D49").Value = 4219.808790577632
$ws.Range("B50").Value = 15882.40260383775
$ws.Range("C50").Value = 663.03957339252
$ws.Range("D50").Value = 4042.576813177311
$ws.Range("B51").Value = 15991.98449512698
$ws.Range("C51").Value = 680.0236443376853
$ws.Range("D51").Value = 3884.867874989115
$ws.Range("B52").Value = 16038.87273219378
$ws.Range("C52").Value = 684.1602477308907
$ws.Range("D52").Value = 3694.313656268508
$ws.Range("B53").Value = 16146.87503836606
$ws.Range("C53").Value = 688.6961366589838
$ws.Range("D53").Value = 3521.044733130797
$ws.Range("B54").Value = 16171.25157442826
$ws.Range("C54").Value = 697.844832067479
$ws.Range("D54").Value = 3353.597435200341
$ws.Range("B55").Value = 16326.80931032308
$ws.Range("C55").Value = 710.5484596831826
$ws.Range("D55").Value = 3241.326560632477
$ws.Range("B56").Value = 16292.45509438363
$ws.Range("C56").Value = 703.7043182293438
$ws.Range("D56").Value = 3028.170349030047
$ws.Range("B57").Value = 16391.6557173154
$ws.Range("C57").Value = 719.6538677611444
$ws.Range("D57").Value = 2851.684331980192
$ws.Range("B58").Value = 16310.10330524439
$ws.Range("C58").Value = 720.8222740452042
$ws.Range("D58").Value = 2725.444665999593
$ws.Range("B59").Value = 16328.71322415624
$ws.Range("C59").Value = 719.7241746379351
$ws.Range("D59").Value = 2589.574228227048
$ws.Range("B60").Value = 16491.15666107844
$ws.Range("C60").Value = 727.074726067095
$ws.Range("D60").Value = 2469.21720801251
$ws.Range("B61").Value = 16498.49445850071
$ws.Range("C61").Value = 725.6516637540846
$ws.Range("D61").Value = 2381.329938304832
$ws.Range("B62").Value = 16517.7241203987
$ws.Range("C62").Value = 728.3680258620219
$ws.Range("D62").Value = 2303.552214656896

Write-Host "Updated trajectory values."
